# Prepare public release: split the single-sheet sample workbook into three
# named sheets with fresh sample data, and trim the stray second data row
# from the first sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet1: rename existing sheet, refresh its sample data -----------------
$ws1 = $wb.ActiveSheet
$ws1.Name = "Sheet1"
$ws1.Range("A1").Value = "Data1"
$ws1.Range("B1").Value = "Value1"
# Drop the old second row ("Test"/"123") entirely so the sheet's used range
# shrinks back down to a single row.
$ws1.Range("A2:B2").ClearContents()

# --- Sheet2: new sheet with its own sample data ------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = "Data2"
$ws2.Range("B1").Value = "Value2"

# --- Sheet3: new, empty sheet -------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

# Restore Sheet1 as the active/selected sheet.
$ws1.Activate()
